$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.803.15"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.316.80"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.51"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.87"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.45"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  +4.16%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "2.680.91"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "2.320.41"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "42.746.92"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.15"
$ws.Range("E19").Value = "  -5.78%  "
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.90"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.26"
$ws.Range("E23").Value = "  +6.22%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.94"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.33"
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.37"
$ws.Range("E28").Value = "  -0.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.30"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.14"
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.26"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.94"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.01"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.46"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0699"
$ws.Range("E36").Value = "  +2.45%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1000"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  +2.89%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.64"
$ws.Range("E42").Value = "  +12.54%  "
$ws.Range("D43").Value = "1.928.05"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.15"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").Value = "2.548.17"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.34"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.12"
$ws.Range("E51").Value = "  +2.59%  "
